$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "Apparel"
    3  = "Aerospace & Defence"
    5  = "Chemicals"
    6  = "Pharmaceuticals"
    8  = "Automotive"
    12 = "Chemicals"
    13 = "Automotive"
    16 = "Logistics"
    17 = "Telecommunication"
    21 = "Construction Materials"
    24 = "Automotive"
    25 = "Pharmaceuticals"
    26 = "Aerospace & Defence"
    28 = "Automotive"
    29 = "Automotive"
    30 = "Biotech"
    31 = "Aerospace & Defence"
    34 = "Medical Technology"
    37 = "Medical Equipment"
    39 = "Automotive"
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
